# changed algorithm to select best sparse matrix
#
# Populate the "general info" sheet with a small Question/Answer table and
# widen the two columns so the new text is readable.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "Question"
$ws.Range("B1").Value = "Answer"

# Data row
$ws.Range("A2").Value = "What is rose-hulman ranking"
$ws.Range("B2").Value = "Rose-hulman is ranked number one"

# Widen the columns to fit the new content.
$ws.Columns.Item(1).ColumnWidth = 29.333333333333332
$ws.Columns.Item(2).ColumnWidth = 51.0

# Leave the selection where it ends up after entering the table (one row
# below the last row of data, in column B).
[void]$ws.Range("B3").Select()
